# Generate Report for Handoff
# Adds two new localization entries (4df4343e..., 8b3ea37d...) ahead of the
# existing e4b246de... entry (status: "Ready for handoff") across the
# Overview, zh-cn and de-de sheets, and grows the three tables accordingly.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

# Clear the old hyperlinks - they get fully rebuilt below in final order.
$wsOv.Hyperlinks.Delete()

# Row 3: now describes 4df4343e... (previously described e4b246de...)
$wsOv.Range("A3").Value = "4df4343e-8da6-4563-b670-4855371533e4.md"
$wsOv.Range("B3").Value = "e2e\4df4343e-8da6-4563-b670-4855371533e4.md"
$wsOv.Range("C3").Value = ".md"
$wsOv.Range("D3").Value = ""
$wsOv.Range("E3").Value = "Ready for handoff"
$wsOv.Range("F3").Value = "Ready for handoff"
$wsOv.Range("G3").Value = "2016-08-24 16:42:26"
$wsOv.Range("G3").NumberFormat = $dateFmt

# Row 4 (new): 8b3ea37d...
$wsOv.Range("A4").Value = "8b3ea37d-c86c-48ba-beda-d897d9c068ee.md"
$wsOv.Range("B4").Value = "e2e\8b3ea37d-c86c-48ba-beda-d897d9c068ee.md"
$wsOv.Range("C4").Value = ".md"
$wsOv.Range("D4").Value = ""
$wsOv.Range("E4").Value = "Ready for handoff"
$wsOv.Range("F4").Value = "Ready for handoff"
$wsOv.Range("G4").Value = "2016-08-24 16:42:26"
$wsOv.Range("G4").NumberFormat = $dateFmt

# Row 5 (new): e4b246de... (moved down from row 3)
$wsOv.Range("A5").Value = "e4b246de-4150-4f3c-9365-8678f5f86480.md"
$wsOv.Range("B5").Value = "e2e\e4b246de-4150-4f3c-9365-8678f5f86480.md"
$wsOv.Range("C5").Value = ".md"
$wsOv.Range("D5").Value = ""
$wsOv.Range("E5").Value = "Ready for handoff"
$wsOv.Range("F5").Value = "Ready for handoff"
$wsOv.Range("G5").Value = "2016-08-24 16:40:59"
$wsOv.Range("G5").NumberFormat = $dateFmt

$wsOv.Hyperlinks.Add($wsOv.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f84fcb9d60bf2544a31ed45ec12383c67d1ebb40/e2e/b0701422-b275-4b16-a8db-2b8cb6dc240f.md", "", "", "e2e\b0701422-b275-4b16-a8db-2b8cb6dc240f.md")
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4df4343e8da64563b6704855371533e4f84fcb9/e2e/4df4343e-8da6-4563-b670-4855371533e4.md", "", "", "e2e\4df4343e-8da6-4563-b670-4855371533e4.md")
$wsOv.Hyperlinks.Add($wsOv.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8b3ea37dc86c48badep897d9c068eef84fcb9d6/e2e/8b3ea37d-c86c-48ba-beda-d897d9c068ee.md", "", "", "e2e\8b3ea37d-c86c-48ba-beda-d897d9c068ee.md")
$wsOv.Hyperlinks.Add($wsOv.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d82bfb72b6218bfc5bfadd5b1057cc1d746044cf/e2e/e4b246de-4150-4f3c-9365-8678f5f86480.md", "", "", "e2e\e4b246de-4150-4f3c-9365-8678f5f86480.md")

$wsOv.ListObjects.Item("Overview").Resize($wsOv.Range("A1:G5"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Delete()

# Row 3: 4df4343e...
$wsZh.Range("A3").Value = "4df4343e-8da6-4563-b670-4855371533e4.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = "4df4343e-8da6-4563-b670-4855371533e4.5d9404d19d94d9be6032f7b861748adb996647db.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-24 16:42:21"
$wsZh.Range("H3").NumberFormat = $dateFmt
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("K3").NumberFormat = $dateFmt
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

# Row 4 (new): 8b3ea37d...
$wsZh.Range("A4").Value = "8b3ea37d-c86c-48ba-beda-d897d9c068ee.md"
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "False"
$wsZh.Range("G4").Value = "8b3ea37d-c86c-48ba-beda-d897d9c068ee.81d04425ecc63c9a3c00575d348a28b9d75c34f4.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-08-24 16:42:21"
$wsZh.Range("H4").NumberFormat = $dateFmt
$wsZh.Range("I4").Value = ""
$wsZh.Range("J4").Value = ""
$wsZh.Range("K4").Value = "0001-01-01 00:00:00"
$wsZh.Range("K4").NumberFormat = $dateFmt
$wsZh.Range("L4").Value = ""
$wsZh.Range("M4").Value = "True"
$wsZh.Range("N4").Value = ""
$wsZh.Range("O4").Value = "False"
$wsZh.Range("P4").Value = ""

# Row 5 (new): e4b246de... (moved down from row 3)
$wsZh.Range("A5").Value = "e4b246de-4150-4f3c-9365-8678f5f86480.md"
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "e2e"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("F5").Value = "False"
$wsZh.Range("G5").Value = "e4b246de-4150-4f3c-9365-8678f5f86480.9dc3f50b3eedf6df04281c7f92c67460c609989c.zh-cn.xlf"
$wsZh.Range("H5").Value = "2016-08-24 16:40:54"
$wsZh.Range("H5").NumberFormat = $dateFmt
$wsZh.Range("I5").Value = ""
$wsZh.Range("J5").Value = ""
$wsZh.Range("K5").Value = "0001-01-01 00:00:00"
$wsZh.Range("K5").NumberFormat = $dateFmt
$wsZh.Range("L5").Value = ""
$wsZh.Range("M5").Value = "True"
$wsZh.Range("N5").Value = ""
$wsZh.Range("O5").Value = "False"
$wsZh.Range("P5").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f84fcb9d60bf2544a31ed45ec12383c67d1ebb40/e2e/b0701422-b275-4b16-a8db-2b8cb6dc240f.md", "", "", "b0701422-b275-4b16-a8db-2b8cb6dc240f.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c975ed7ac8133219e39e0d178e763c0ecd5f2e40/e2e/b0701422-b275-4b16-a8db-2b8cb6dc240f.md", "", "", "b0701422-b275-4b16-a8db-2b8cb6dc240f.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4df4343e8da64563b6704855371533e4f84fcb9/e2e/4df4343e-8da6-4563-b670-4855371533e4.md", "", "", "4df4343e-8da6-4563-b670-4855371533e4.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8b3ea37dc86c48badep897d9c068eef84fcb9d6/e2e/8b3ea37d-c86c-48ba-beda-d897d9c068ee.md", "", "", "8b3ea37d-c86c-48ba-beda-d897d9c068ee.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d82bfb72b6218bfc5bfadd5b1057cc1d746044cf/e2e/e4b246de-4150-4f3c-9365-8678f5f86480.md", "", "", "e4b246de-4150-4f3c-9365-8678f5f86480.md")

$wsZh.ListObjects.Item("zh-cn").Resize($wsZh.Range("A1:P5"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Delete()

# Row 3: 4df4343e...
$wsDe.Range("A3").Value = "4df4343e-8da6-4563-b670-4855371533e4.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = "4df4343e-8da6-4563-b670-4855371533e4.5d9404d19d94d9be6032f7b861748adb996647db.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-24 16:42:26"
$wsDe.Range("H3").NumberFormat = $dateFmt
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("K3").NumberFormat = $dateFmt
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

# Row 4 (new): 8b3ea37d...
$wsDe.Range("A4").Value = "8b3ea37d-c86c-48ba-beda-d897d9c068ee.md"
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "False"
$wsDe.Range("G4").Value = "8b3ea37d-c86c-48ba-beda-d897d9c068ee.81d04425ecc63c9a3c00575d348a28b9d75c34f4.de-de.xlf"
$wsDe.Range("H4").Value = "2016-08-24 16:42:26"
$wsDe.Range("H4").NumberFormat = $dateFmt
$wsDe.Range("I4").Value = ""
$wsDe.Range("J4").Value = ""
$wsDe.Range("K4").Value = "0001-01-01 00:00:00"
$wsDe.Range("K4").NumberFormat = $dateFmt
$wsDe.Range("L4").Value = ""
$wsDe.Range("M4").Value = "True"
$wsDe.Range("N4").Value = ""
$wsDe.Range("O4").Value = "False"
$wsDe.Range("P4").Value = ""

# Row 5 (new): e4b246de... (moved down from row 3)
$wsDe.Range("A5").Value = "e4b246de-4150-4f3c-9365-8678f5f86480.md"
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "e2e"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("F5").Value = "False"
$wsDe.Range("G5").Value = "e4b246de-4150-4f3c-9365-8678f5f86480.9dc3f50b3eedf6df04281c7f92c67460c609989c.de-de.xlf"
$wsDe.Range("H5").Value = "2016-08-24 16:40:59"
$wsDe.Range("H5").NumberFormat = $dateFmt
$wsDe.Range("I5").Value = ""
$wsDe.Range("J5").Value = ""
$wsDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDe.Range("K5").NumberFormat = $dateFmt
$wsDe.Range("L5").Value = ""
$wsDe.Range("M5").Value = "True"
$wsDe.Range("N5").Value = ""
$wsDe.Range("O5").Value = "False"
$wsDe.Range("P5").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f84fcb9d60bf2544a31ed45ec12383c67d1ebb40/e2e/b0701422-b275-4b16-a8db-2b8cb6dc240f.md", "", "", "b0701422-b275-4b16-a8db-2b8cb6dc240f.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1df39bc1db6bd17a2f42cec1e878130c7839252e/e2e/b0701422-b275-4b16-a8db-2b8cb6dc240f.md", "", "", "b0701422-b275-4b16-a8db-2b8cb6dc240f.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4df4343e8da64563b6704855371533e4f84fcb9/e2e/4df4343e-8da6-4563-b670-4855371533e4.md", "", "", "4df4343e-8da6-4563-b670-4855371533e4.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8b3ea37dc86c48badep897d9c068eef84fcb9d6/e2e/8b3ea37d-c86c-48ba-beda-d897d9c068ee.md", "", "", "8b3ea37d-c86c-48ba-beda-d897d9c068ee.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d82bfb72b6218bfc5bfadd5b1057cc1d746044cf/e2e/e4b246de-4150-4f3c-9365-8678f5f86480.md", "", "", "e4b246de-4150-4f3c-9365-8678f5f86480.md")

$wsDe.ListObjects.Item("de-de").Resize($wsDe.Range("A1:P5"))
